{"js": "const pairs = [\n  [\"51\u00d740=\", \"72\u00d750=\"],\n  [\"48\u00d790=\", \"43\u00d757=\"],\n  [\"17\u00d734=\", \"42\u00d746=\"],\n  [\"81\u00d766=\", \"68\u00d712=\"],\n  [\"97\u00d742=\", \"93\u00d736=\"],\n  [\"97\u00d793=\", \"46\u00d789=\"],\n  [\"85\u00d772=\", \"97\u00d759=\"],\n  [\"90\u00d751=\", \"63\u00d719=\"],\n  [\"58\u00d741=\", \"35\u00d768=\"],\n  [\"39\u00d770=\", \"24\u00d753=\"],\n  [\"86\u00d791=\", \"27\u00d773=\"],\n  [\"74\u00d740=\", \"20\u00d797=\"],\n  [\"38\u00d757=\", \"16\u00d760=\"],\n  [\"18\u00d742=\", \"44\u00d760=\"],\n  [\"54\u00d757=\", \"47\u00d793=\"],\n  [\"47\u00d778=\", \"44\u00d748=\"],\n  [\"21\u00d773=\", \"69\u00d792=\"],\n  [\"53\u00d739=\", \"16\u00d722=\"],\n  [\"54\u00d794=\", \"39\u00d721=\"],\n  [\"62\u00d740=\", \"57\u00d750=\"],\n  [\"41\u00d771=\", \"78\u00d712=\"],\n  [\"42\u00d750=\", \"91\u00d777=\"],\n  [\"44\u00d725=\", \"17\u00d730=\"],\n  [\"43\u00d763=\", \"46\u00d741=\"],\n  [\"86\u00d720=\", \"44\u00d795=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"51\u00d740=\", \"72\u00d750=\"),\n    @(\"48\u00d790=\", \"43\u00d757=\"),\n    @(\"17\u00d734=\", \"42\u00d746=\"),\n    @(\"81\u00d766=\", \"68\u00d712=\"),\n    @(\"97\u00d742=\", \"93\u00d736=\"),\n    @(\"97\u00d793=\", \"46\u00d789=\"),\n    @(\"85\u00d772=\", \"97\u00d759=\"),\n    @(\"90\u00d751=\", \"63\u00d719=\"),\n    @(\"58\u00d741=\", \"35\u00d768=\"),\n    @(\"39\u00d770=\", \"24\u00d753=\"),\n    @(\"86\u00d791=\", \"27\u00d773=\"),\n    @(\"74\u00d740=\", \"20\u00d797=\"),\n    @(\"38\u00d757=\", \"16\u00d760=\"),\n    @(\"18\u00d742=\", \"44\u00d760=\"),\n    @(\"54\u00d757=\", \"47\u00d793=\"),\n    @(\"47\u00d778=\", \"44\u00d748=\"),\n    @(\"21\u00d773=\", \"69\u00d792=\"),\n    @(\"53\u00d739=\", \"16\u00d722=\"),\n    @(\"54\u00d794=\", \"39\u00d721=\"),\n    @(\"62\u00d740=\", \"57\u00d750=\"),\n    @(\"41\u00d771=\", \"78\u00d712=\"),\n    @(\"42\u00d750=\", \"91\u00d777=\"),\n    @(\"44\u00d725=\", \"17\u00d730=\"),\n    @(\"43\u00d763=\", \"46\u00d741=\"),\n    @(\"86\u00d720=\", \"44\u00d795=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
